$wb = $excel.ActiveWorkbook


# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 834.1177
$ws.Range("J28").Value = 1801.2
$ws.Range("L28").Value = 1801.2
$ws.Range("N28").Value = -2771.2
$ws.Range("H33").Value = 325667.62
$ws.Range("I33").Value = 528.1539
$ws.Range("J33").Value = 709923.4
$ws.Range("K33").Value = 528.1539
$ws.Range("L33").Value = 709923.4
$ws.Range("M33").Value = -299.1539
$ws.Range("N33").Value = -710381.4
$ws.Range("H41").Value = 908.3333
$ws.Range("I41").Value = 1142.0667
$ws.Range("J41").Value = 674.6
$ws.Range("K41").Value = 1142.0667
$ws.Range("L41").Value = 674.6
$ws.Range("M41").Value = -702.0667000000001
$ws.Range("N41").Value = -1554.6
$ws.Range("H80").Value = 48809.047
$ws.Range("I80").Value = 1158.7
$ws.Range("J80").Value = 92127.55
$ws.Range("K80").Value = 3476.1
$ws.Range("L80").Value = 276382.65
$ws.Range("M80").Value = -2478.1
$ws.Range("N80").Value = -278378.65
$ws.Range("H83").Value = 48809.047
$ws.Range("I83").Value = 1158.7
$ws.Range("J83").Value = 92127.55
$ws.Range("K83").Value = 10428.3
$ws.Range("L83").Value = 829147.9500000001
$ws.Range("M83").Value = -5436.300000000001
$ws.Range("N83").Value = -839131.9500000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3114.93
$ws.Range("I32").Value = 2798.897
$ws.Range("J32").Value = 13333.333
$ws.Range("K32").Value = 2798.897
$ws.Range("L32").Value = 13333.333
$ws.Range("M32").Value = -2511.897
$ws.Range("N32").Value = -13907.333
$ws.Range("H61").Value = 2920.077
$ws.Range("I61").Value = 2240.0
$ws.Range("J61").Value = 3082.0
$ws.Range("K61").Value = 2240.0
$ws.Range("L61").Value = 3082.0
$ws.Range("M61").Value = -2028.0
$ws.Range("N61").Value = -3506.0
$ws.Range("H122").Value = 1752.2
$ws.Range("I122").Value = 1674.3529
$ws.Range("J122").Value = 2193.3333
$ws.Range("K122").Value = 5023.0587
$ws.Range("L122").Value = 6579.999899999999
$ws.Range("M122").Value = -2573.0587
$ws.Range("N122").Value = -11479.9999
$ws.Range("H136").Value = 2920.077
$ws.Range("I136").Value = 2240.0
$ws.Range("J136").Value = 3082.0
$ws.Range("K136").Value = 6720.0
$ws.Range("L136").Value = 9246.0
$ws.Range("M136").Value = -4170.0
$ws.Range("N136").Value = -14346.0

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1500.762
$ws.Range("I64").Value = 720.4
$ws.Range("J64").Value = 1744.625
$ws.Range("K64").Value = 720.4
$ws.Range("L64").Value = 1744.625
$ws.Range("M64").Value = -495.4
$ws.Range("N64").Value = -2194.625
$ws.Range("H67").Value = 1500.762
$ws.Range("I67").Value = 720.4
$ws.Range("J67").Value = 1744.625
$ws.Range("K67").Value = 720.4
$ws.Range("L67").Value = 1744.625
$ws.Range("M67").Value = 59.60000000000002
$ws.Range("N67").Value = -3304.625
$ws.Range("H86").Value = 102227.82
$ws.Range("I86").Value = 112290.6
$ws.Range("J86").Value = 1600.0
$ws.Range("K86").Value = 112290.6
$ws.Range("L86").Value = 1600.0
$ws.Range("M86").Value = -111167.6
$ws.Range("N86").Value = -3846.0
$ws.Range("H89").Value = 102227.82
$ws.Range("I89").Value = 112290.6
$ws.Range("J89").Value = 1600.0
$ws.Range("K89").Value = 561453.0
$ws.Range("L89").Value = 8000.0
$ws.Range("M89").Value = -555837.0
$ws.Range("N89").Value = -19232.0
$ws.Range("H94").Value = 729.7
$ws.Range("I94").Value = 723.17645
$ws.Range("J94").Value = 766.6667
$ws.Range("K94").Value = 723.17645
$ws.Range("L94").Value = 766.6667
$ws.Range("M94").Value = -272.17645
$ws.Range("N94").Value = -1668.6667
$ws.Range("H130").Value = 32554.264
$ws.Range("J130").Value = 32554.264
$ws.Range("L130").Value = 32554.264
$ws.Range("N130").Value = -42594.264

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2000.0
$ws.Range("I16").Value = 0.0
$ws.Range("J16").Value = 2000.0
$ws.Range("K16").Value = 0.0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -2574.0
$ws.Range("H86").Value = 3685.739
$ws.Range("I86").Value = 3309.625
$ws.Range("J86").Value = 3886.3333
$ws.Range("K86").Value = 3309.625
$ws.Range("L86").Value = 3886.3333
$ws.Range("M86").Value = -2186.625
$ws.Range("N86").Value = -6132.3333
$ws.Range("H89").Value = 3685.739
$ws.Range("I89").Value = 3309.625
$ws.Range("J89").Value = 3886.3333
$ws.Range("K89").Value = 16548.125
$ws.Range("L89").Value = 19431.6665
$ws.Range("M89").Value = -10932.125
$ws.Range("N89").Value = -30663.6665
$ws.Range("H113").Value = 2000.0
$ws.Range("I113").Value = 0.0
$ws.Range("J113").Value = 2000.0
$ws.Range("K113").Value = 0.0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6340.0
$ws.Range("H122").Value = 1200.0
$ws.Range("I122").Value = 0.0
$ws.Range("J122").Value = 1200.0
$ws.Range("K122").Value = 0.0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -8500.0

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 658.2222
$ws.Range("I122").Value = 848.5
$ws.Range("J122").Value = 634.4375
$ws.Range("K122").Value = 7636.5
$ws.Range("L122").Value = 5709.9375
$ws.Range("M122").Value = -5186.5
$ws.Range("N122").Value = -10609.9375
$ws.Range("H131").Value = 836.6869
$ws.Range("I131").Value = 565.0
$ws.Range("J131").Value = 848.12634
$ws.Range("K131").Value = 1695.0
$ws.Range("L131").Value = 2544.37902
$ws.Range("M131").Value = 3345.0
$ws.Range("N131").Value = -12624.37902

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 58886324.0
$ws.Range("I80").Value = 100104910.0
$ws.Range("J80").Value = 2628.8572
$ws.Range("K80").Value = 100104910.0
$ws.Range("L80").Value = 2628.8572
$ws.Range("M80").Value = -100103912.0
$ws.Range("N80").Value = -4624.8572
$ws.Range("H83").Value = 58886324.0
$ws.Range("I83").Value = 100104910.0
$ws.Range("J83").Value = 2628.8572
$ws.Range("K83").Value = 500524550.0
$ws.Range("L83").Value = 13144.286
$ws.Range("M83").Value = -500519558.0
$ws.Range("N83").Value = -23128.286
$ws.Range("H97").Value = 45456840.0
$ws.Range("I97").Value = 62502456.0
$ws.Range("J97").Value = 1868.3334
$ws.Range("K97").Value = 62502456.0
$ws.Range("L97").Value = 1868.3334
$ws.Range("M97").Value = -62501960.0
$ws.Range("N97").Value = -2860.3334
$ws.Range("H102").Value = 403977.72
$ws.Range("I102").Value = 2308.7778
$ws.Range("J102").Value = 1006481.2
$ws.Range("K102").Value = 2308.7778
$ws.Range("L102").Value = 1006481.2
$ws.Range("M102").Value = -686.7777999999998
$ws.Range("N102").Value = -1009725.2
$ws.Range("H134").Value = 45500.0
$ws.Range("J134").Value = 45500.0
$ws.Range("L134").Value = 136500.0
$ws.Range("N134").Value = -141570.0

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1956.8572
$ws.Range("I22").Value = 2673.5
$ws.Range("J22").Value = 1670.2
$ws.Range("K22").Value = 2673.5
$ws.Range("L22").Value = 1670.2
$ws.Range("M22").Value = -2378.5
$ws.Range("N22").Value = -2260.2
$ws.Range("H27").Value = 1956.8572
$ws.Range("I27").Value = 2673.5
$ws.Range("J27").Value = 1670.2
$ws.Range("K27").Value = 2673.5
$ws.Range("L27").Value = 1670.2
$ws.Range("M27").Value = -2566.5
$ws.Range("N27").Value = -1884.2

Write-Output "Applied 191 cell edits across 7 sheets"